$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.508.65'
$ws.Range('E2').Value = '  -4.44%  '
$ws.Range('D3').Value = '2.957.76'
$ws.Range('E3').Value = '  -6.41%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.67'
$ws.Range('E5').Value = '  -5.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.08'
$ws.Range('E6').Value = '  -7.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('D9').Value = '2.964.73'
$ws.Range('E9').Value = '  -6.19%  '
$ws.Range('E10').Value = '  -3.41%  '
$ws.Range('E11').Value = '  -7.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.368'
$ws.Range('E12').Value = '  -3.77%  '
$ws.Range('D13').Value = '3.476.94'
$ws.Range('E13').Value = '  -6.49%  '
$ws.Range('E14').Value = '  -2.89%  '
$ws.Range('D15').Value = '61.578.66'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.70'
$ws.Range('E16').Value = '  -5.84%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000147'
$ws.Range('E17').Value = '  -4.96%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.960.51'
$ws.Range('E18').Value = '  -6.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.16'
$ws.Range('E19').Value = '  -1.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '382.40'
$ws.Range('E20').Value = '  -5.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.98'
$ws.Range('E21').Value = '  -5.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.67'
$ws.Range('E22').Value = '  -6.23%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.18'
$ws.Range('E24').Value = '  -5.15%  '
$ws.Range('E25').Value = '  -3.13%  '
$ws.Range('D26').Value = '3.083.67'
$ws.Range('E26').Value = '  -6.64%  '
$ws.Range('E27').Value = '  -3.98%  '
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').Value = '0.0₃0935'
$ws.Range('E29').Value = '  -8.29%  '
$ws.Range('E30').Value = '  -5.25%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -5.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.43'
$ws.Range('E33').Value = '  -3.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '159.29'
$ws.Range('E34').Value = '  +1.67%  '
$ws.Range('E35').Value = '  -3.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.95'
$ws.Range('E36').Value = '  -5.47%  '
$ws.Range('E37').Value = '  -4.60%  '
$ws.Range('E38').Value = '  -4.73%  '
$ws.Range('E39').Value = '  -7.36%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.94'
$ws.Range('E40').Value = '  -3.56%  '
$ws.Range('D41').Value = '2.410.65'
$ws.Range('E41').Value = '  -9.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.21'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.14'
$ws.Range('E43').Value = '  -7.10%  '
$ws.Range('E44').Value = '  -4.73%  '
$ws.Range('E45').Value = '  -3.29%  '
$ws.Range('E46').Value = '  -0.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0246'
$ws.Range('E47').Value = '  -3.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.97'
$ws.Range('E48').Value = '  -8.80%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0958'
$ws.Range('E49').Value = '  -2.37%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '268.13'
$ws.Range('E50').Value = '  -7.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.77'
$ws.Range('E51').Value = '  -7.04%  '
